$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# PART A: Reorganize the "2018/09/11" block.
# ---------------------------------------------------------------------------
# 1) Duplicate paragraph 1 (date header "2018/09/11: doannd2") and paste the
#    copy right after paragraph 2 ("Thêm GUI..."), i.e. before paragraph 3.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Copy()
$target1 = $d.Range($d.Paragraphs.Item(2).Range.End, $d.Paragraphs.Item(2).Range.End)
$target1.Paste()

# 2) Duplicate (original) paragraph 2 ("Thêm GUI...") and paste the copy
#    right after the newly-inserted date header (now paragraph 3), i.e.
#    before the old paragraph 3 ("Cửa sổ ...").
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Copy()
$target2 = $d.Range($d.Paragraphs.Item(3).Range.End, $d.Paragraphs.Item(3).Range.End)
$target2.Paste()

# 3) Edit the ORIGINAL paragraph 1: change the date from 2018/09/11 to
#    2018/09/12. Target only the single character inside the middle run
#    (which holds "1") so the run structure ("2018/09/1" | "1" | ": doannd2")
#    stays intact - only its text content changes.
$p1b = $d.Paragraphs.Item(1)
$dateDigitStart = $p1b.Range.Start + 9   # "2018/09/1" is 9 characters
$dateDigitRange = $d.Range($dateDigitStart, $dateDigitStart + 1)
$dateDigitRange.Text = "2"

# 4) Edit the ORIGINAL paragraph 2: replace its text with the new changelog
#    entry, and move the "_GoBack" bookmark from paragraph 4 (now shifted)
#    to the end of this paragraph.
$p2b = $d.Paragraphs.Item(2)
$p2bTextRange = $p2b.Range.Duplicate
$p2bTextRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$p2bTextRange.Text = "Bỏ schedule cho thu hoạch, thay bằng setTimeOut"

# 5) Move the "_GoBack" bookmark from its current paragraph (the "Cache
#    user name..." item) to the end of paragraph 2 ("Bỏ schedule...").
#    A collapsed bookmark placed exactly at (paragraph.End - 1) is mishandled
#    by this runtime, so we work around it: insert a temporary placeholder
#    character, wrap a non-collapsed bookmark around it, then delete the
#    placeholder through the bookmark's own range (which leaves behind a
#    correctly-positioned collapsed bookmark).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$p2c = $d.Paragraphs.Item(2)
$insPos = $p2c.Range.End - 1
$ph = $d.Range($insPos, $insPos)
$ph.InsertAfter("X")
$bmRange = $d.Range($insPos, $insPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Bookmarks.Item("_GoBack").Range.Text = ""

Write-Host "Paragraphs now:" $d.Paragraphs.Count
for ($i=1; $i -le 7; $i++) {
    Write-Host "Para $i : [$($d.Paragraphs.Item($i).Range.Text)]"
}
